$wb = $excel.ActiveWorkbook

# Add the new "pilotos" worksheet right after the existing "familia" sheet
$familia = $wb.Worksheets.Item("familia")
$pilotos = $wb.Worksheets.Add($null, $familia)
$pilotos.Name = "pilotos"

# Headers
$pilotos.Range("A1").Value = "nombre"
$pilotos.Range("B1").Value = "apellido"
$pilotos.Range("C1").Value = "numero"
$pilotos.Range("D1").Value = "escuderia"

# Row 2
$pilotos.Range("A2").Value = "Charles"
$pilotos.Range("B2").Value = "Leclerc"
$pilotos.Range("C2").Value = 16
$pilotos.Range("D2").Value = "ferrari"

# Row 3
$pilotos.Range("A3").Value = "Max"
$pilotos.Range("B3").Value = "Verstappen"
$pilotos.Range("C3").Value = 1
$pilotos.Range("D3").Value = "red bull"

# Row 4
$pilotos.Range("A4").Value = "Lewis"
$pilotos.Range("B4").Value = "Hamilton"
$pilotos.Range("C4").Value = 44
$pilotos.Range("D4").Value = "mercedes"

# Select the whole data range on "familia" (no cell in particular highlighted)
$familia.Range("A1:C4").Select()

# Activate the new "pilotos" sheet and position the cursor just below its data
$pilotos.Activate()
$pilotos.Range("D5").Select()
